# lecture: keyword_id -> question_keywords: foreign key로(lecture_id 대신 keyword_id)
#
# 1. LECTURE_KEYWORDS.id  -> LECTURE_KEYWORDS.keyword_id   (column name only, row 22)
# 2. QUESTION_KEYWORDS.keyword (VARCHAR(255)) -> QUESTION_KEYWORDS.keyword_id (INT(11))  (row 37)
# 3. QUESTION_KEYWORDS.lecture_id row (old row 38) is removed entirely - the table
#    now references LECTURE_KEYWORDS.keyword_id as its foreign key instead of lecture_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) LECTURE_KEYWORDS primary key column rename: id -> keyword_id
$ws.Range("C22").Value = "keyword_id"

# 2) QUESTION_KEYWORDS: keyword/키워드/VARCHAR(255) becomes keyword_id/키워드 ID/INT(11)
$ws.Range("C37").Value = "keyword_id"
$ws.Range("D37").Value = "키워드 ID"
$ws.Range("E37").Value = "INT(11)"

# 3) Remove the old lecture_id row (row 38) entirely - rows below shift up by one
$ws.Rows(38).Delete()

# Restore the selection to match the post-edit view
$ws.Range("D37").Select() | Out-Null
